# Switch licence from BY-NC to BY-SA
#
# The instructor-guide footer states that the work is licensed under the
# Creative Commons "BY-NC 4.0" licence and links to its legal-code page.
# This script updates both the human-readable licence name in the body
# text and the Creative Commons hyperlink (its display text and its
# target URL) so that everything instead refers to "BY-SA 4.0".

$d = $word.ActiveDocument

# 1. Body text: "... is licensed under CC BY-NC 4.0. To view ..."
#             -> "... is licensed under CC BY-SA 4.0. To view ..."
$null = $d.Content.Find.Execute("CC BY-NC 4.0", $false, $false, $false, $false, `
                                 $false, $true, 1, $false, "CC BY-SA 4.0", 2)

# 2. Hyperlink that points at the CC BY-NC licence -> point it at BY-SA,
#    and update the displayed URL text to match.
foreach ($hl in $d.Hyperlinks) {
    if ($hl.Address -like "*creativecommons.org/licenses/by-nc/4.0*") {
        $hl.TextToDisplay = "https://creativecommons.org/licenses/by-sa/4.0"
        $hl.Address = "https://creativecommons.org/licenses/by-sa/4.0"
    }
}
